$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date bumped by one day (45183 -> 45184) for every
# data row (2 through 89).
for ($r = 2; $r -le 89; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}

# Row 2 (A 33491-2023 / HOFORS) link-formula touch-up: add the friendly
# display-name second argument to each HYPERLINK() call, and turn the
# Y2 cell (previously stored as plain inline text) into a real formula.
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/artfynd/A 33491-2023.xlsx, "A 33491-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/kartor/A 33491-2023.png", "A 33491-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/klagomål/A 33491-2023.docx", "A 33491-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/klagomålsmail/A 33491-2023.docx", "A 33491-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/tillsyn/A 33491-2023.docx", "A 33491-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_HOFORS/tillsynsmail/A 33491-2023.docx", "A 33491-2023")'
